$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.928.75'
$ws.Range("D3").Value = '1.552.21'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("E4").Value = '  -0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.488'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.96'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.247'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0856'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").Value = '1.773.40'
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").Value = '1.558.22'
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").Value = '26.917.17'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").Value = '0.0₃0711'
$ws.Range("E18").Value = '  +3.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.31%  '
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("E32").Value = '  -0.36%  '
$ws.Range("E33").Value = '  +3.49%  '
$ws.Range("D34").Value = '1.411.03'
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0165'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.525'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("E40").Value = '  -0.50%  '
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("E42").Value = '  +3.18%  '
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").Value = '1.686.99'
$ws.Range("E47").Value = '  -0.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.30'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0520'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.72%  '
$ws.Range("E50").Value = '  +2.75%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0959'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.06%  '
